$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------
# Sheet1 (DeliverynoteHeader)
# ---------------------------------------------------------------

# B2: date value becomes a plain text string, left aligned ("text" numFmt)
$ws1.Range("B2").Value = "16-12-2025"
$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").HorizontalAlignment = -4131

# E1:I1 new bold header cells (Tax, Customer Notes, Terms And Conditions, Save As, Price List)
# Build the desired format once on a helper cell, then copy the format over so the
# engine reuses a single new cell style instead of emitting extra ones.
$ws1.Range("Z1").Font.Bold = $true
$ws1.Range("Z1").WrapText = $false
$ws1.Range("Z1").VerticalAlignment = -4107
$ws1.Range("Z1").Copy()
$ws1.Range("E1:I1").PasteSpecial(-4122)

$ws1.Range("E1").Value = "Tax"
$ws1.Range("F1").Value = "Customer Notes"
$ws1.Range("G1").Value = "Terms And Conditions"
$ws1.Range("H1").Value = "Save As"
$ws1.Range("I1").Value = "Price List"
$ws1.Range("Z1").Clear()

# E2:G2 new text cells (Inclusive, notex, termsx) reuse the existing "text" style
$ws1.Range("Z2").NumberFormat = "@"
$ws1.Range("Z2").HorizontalAlignment = -4131
$ws1.Range("Z2").Copy()
$ws1.Range("E2:G2").PasteSpecial(-4122)

$ws1.Range("E2").Value = "Inclusive"
$ws1.Range("F2").Value = "notex"
$ws1.Range("G2").Value = "termsx"
$ws1.Range("Z2").Clear()

# I2 must be written before H2 so the shared-string table keeps the original order
$ws1.Range("I2").Value = "special price"

# H2 uses a plain (non-bold) Calibri font
$ws1.Range("Z3").Font.Name = "Calibri"
$ws1.Range("Z3").Font.Size = 11
$ws1.Range("Z3").Copy()
$ws1.Range("H2").PasteSpecial(-4122)
$ws1.Range("H2").Value = "SAVE AS DRAFT"
$ws1.Range("Z3").Clear()

# ---------------------------------------------------------------
# Sheet2 (deliverynoteItems)
# ---------------------------------------------------------------

$ws2.Range("D1").Value = "Discount Type"
$ws2.Range("D1").Font.Bold = $true
$ws2.Range("E1").Value = "Discount"
$ws2.Range("E1").Font.Bold = $true

$ws2.Range("D2").Value = "%"
$ws2.Range("E2").Value = 10

$ws2.Range("D3").Value = "amount"
$ws2.Range("E3").Value = 5

$ws2.Range("B3").Value = "Banana"
$ws2.Range("A3").Value = "subinm"
$ws2.Range("C3").Value = 4

# ---------------------------------------------------------------
# View state: active sheet becomes DeliverynoteHeader (first sheet),
# with new selections on both sheets.
# ---------------------------------------------------------------
$ws2.Range("F6").Select()
$ws1.Activate()
$ws1.Range("A4").Select()
